# Update "想去人数" (interested-count) figures in column F for the rows
# whose underlying bilibili event page was refreshed with newer counts.
# The workbook carries the same event list twice: once on sheet "展览"
# and once (merged with a couple of extra "演出" rows) on sheet "全部类型",
# which shifts a few row numbers by +1 starting at the "OCG国潮动漫游戏嘉年华" row.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 14614   # 昆山·第十二届理想乡动漫游戏展
$ws.Range("F5").Value  = 17718   # 苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区
$ws.Range("F15").Value = 44      # 苏州·动漫游戏嘉年华cv见面会
$ws.Range("F17").Value = 153     # 苏州·苏州湾动漫游戏嘉年华
$ws.Range("F19").Value = 1335    # 常熟·CDW·动漫展03
$ws.Range("F24").Value = 7301    # 苏州·OCG国潮动漫游戏嘉年华
$ws.Range("F30").Value = 5867    # 【会员购严选】苏州·Come in joy动漫国潮文化节
$ws.Range("F32").Value = 47      # 苏州·明日方舟ONLY#2024~佑桑柔
$ws.Range("F36").Value = 5102    # 苏州·萤火国潮文化节动漫品牌博览会

# --- Sheet "全部类型" (same events, rows offset by the extra 演出 rows) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 14614   # 昆山·第十二届理想乡动漫游戏展
$ws.Range("F5").Value  = 17718   # 苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区
$ws.Range("F15").Value = 44      # 苏州·动漫游戏嘉年华cv见面会
$ws.Range("F17").Value = 153     # 苏州·苏州湾动漫游戏嘉年华
$ws.Range("F19").Value = 1335    # 常熟·CDW·动漫展03
$ws.Range("F25").Value = 7301    # 苏州·OCG国潮动漫游戏嘉年华
$ws.Range("F32").Value = 5867    # 【会员购严选】苏州·Come in joy动漫国潮文化节
$ws.Range("F34").Value = 47      # 苏州·明日方舟ONLY#2024~佑桑柔
$ws.Range("F38").Value = 5102    # 苏州·萤火国潮文化节动漫品牌博览会
